# Applies the OOXML diff to BPA数据提取.docx:
#  1. Removes the _GoBack bookmark from the first paragraph.
#  2. Adds <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> to the <w:pPr> of the
#     paragraphs in the last row of table 1, and in rows 3-6 of table 4
#     (creating a <w:pPr> for the cells that don't have one).
#  3. Re-adds the _GoBack bookmark at the end of the paragraph in table 4,
#     row 5 ("节点2名称" / "BPA中联络线末端节点名").
#  4. Gives the trailing empty paragraph after table 4 a <w:pPr> with the
#     same eastAsia hint.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------
# 1) Strip the _GoBack bookmark from the document's first paragraph.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.InsertXML("<w:p $wNs><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>BPA</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>提取数据格式</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>：</w:t></w:r></w:p>")

Write-Host "step1 done"

# ---------------------------------------------------------------------
# 2) Table 1 ("分区编号" table), last row: give every cell paragraph a
#    <w:pPr>/<w:rPr><w:rFonts w:hint="eastAsia"/> (creating the <w:pPr>
#    for the 4th cell, which previously had none).
# ---------------------------------------------------------------------
$t1 = $d.Tables(1)
$row = $t1.Rows($t1.Rows.Count)

$row.Cells(1).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>3</w:t></w:r></w:p>")

$row.Cells(2).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>省网编号</w:t></w:r></w:p>")

$row.Cells(3).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>int</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>")

$row.Cells(4).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>分区所在省网编号</w:t></w:r></w:p>")

Write-Host "step2 done"

# ---------------------------------------------------------------------
# 3) Table 4 ("联络线" table), rows 4-7 (row numbers "3","4","5","6" in
#    the printed 列号 column) each get the same eastAsia hint treatment
#    in every cell paragraph. Row 6 ("5" / 节点2名称 / char(8) /
#    "BPA中联络线末端节点名") also gets the _GoBack bookmark re-added at
#    the end of its last cell's paragraph.
# ---------------------------------------------------------------------
$t4 = $d.Tables(4)

# --- row "3" (table row 4): 3 / 省网2编号 / int / 联络线连接的省网 ---
$row = $t4.Rows(4)
$row.Cells(1).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>3</w:t></w:r></w:p>")
$row.Cells(2).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>省网</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>编号</w:t></w:r></w:p>")
$row.Cells(3).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:proofErr w:type=`"spellStart`"/><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>int</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>")
$row.Cells(4).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>联络线连接的省网</w:t></w:r></w:p>")

# --- row "4" (table row 5): 4 / 节点1名称 / char(8) / BPA中联络线首端节点名 ---
$row = $t4.Rows(5)
$row.Cells(1).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>4</w:t></w:r></w:p>")
$row.Cells(2).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>节点</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>名称</w:t></w:r></w:p>")
$row.Cells(3).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>char</w:t></w:r><w:r><w:t>(8)</w:t></w:r></w:p>")
$row.Cells(4).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>BPA</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>中联络线首端节点名</w:t></w:r></w:p>")

# --- row "5" (table row 6): 5 / 节点2名称 / char(8) / BPA中联络线末端节点名 (+ bookmark) ---
$row = $t4.Rows(6)
$row.Cells(1).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>5</w:t></w:r></w:p>")
$row.Cells(2).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>节点</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>名称</w:t></w:r></w:p>")
$row.Cells(3).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>char(8)</w:t></w:r></w:p>")
$row.Cells(4).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>BPA</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>中联络线末端节点名</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>")

# --- row "6" (table row 7): 6 / 基准电压等级 / float / BPA中联络线电压等级 ---
$row = $t4.Rows(7)
$row.Cells(1).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>6</w:t></w:r></w:p>")
$row.Cells(2).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>基准电压等级</w:t></w:r></w:p>")
$row.Cells(3).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>float</w:t></w:r></w:p>")
$row.Cells(4).Range.Paragraphs(1).Range.InsertXML(
  "<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>BPA</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>中联络线电压等级</w:t></w:r></w:p>")

Write-Host "step3 done"

# ---------------------------------------------------------------------
# 4) Trailing empty paragraph after table 4 gets the same eastAsia hint
#    in an (until now absent) <w:pPr>.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$lastPara.InsertXML("<w:p $wNs><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>")

Write-Host "step4 done"
